# Generate Report for Handoff
# - Updates the Priority column ("low" -> "ht") for the four in-flight rows
#   (599551f7…, 90377c4c…, 32dd95cf…, bc536516…) on both locale sheets.
# - Refreshes the "Latest Handoff Datetime" for those same rows on each
#   locale sheet, and the corresponding "Latest HO Xliff Generate Date" on
#   the Overview sheet (de-de handoff datetime is shared with the Overview
#   generate-date column).

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsOverview = $wb.Worksheets.Item("Overview")

# --- zh-cn sheet: Priority (E) and Latest Handoff Datetime (H) for rows 4-7 ---
foreach ($row in 4..7) {
    $wsZhCn.Range("E$row").Value = "ht"
    $wsZhCn.Range("H$row").Value = "2016-09-06 06:04:29"
}

# --- de-de sheet: Priority (E) and Latest Handoff Datetime (H) for rows 4-7 ---
foreach ($row in 4..7) {
    $wsDeDe.Range("E$row").Value = "ht"
    $wsDeDe.Range("H$row").Value = "2016-09-06 06:04:39"
}

# --- Overview sheet: Latest HO Xliff Generate Date (G) for rows 4-7 ---
foreach ($row in 4..7) {
    $wsOverview.Range("G$row").Value = "2016-09-06 06:04:39"
}
